$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency data (price / volume / coin name / link swaps).
# Values are set with a leading apostrophe so Excel stores them as literal text
# (matching the original inlineStr cells) instead of auto-converting numeric-looking
# strings (e.g. '94.344.20', '0.0000201') into numbers, then the style is reset to
# 'Normal' to strip the quote-prefix formatting Excel applies automatically.

$ws.Range("D2").Value = "'94.344.20"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +2.85%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.109.55"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.10%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.24%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'238.14"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -1.73%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'614.68"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -0.10%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +2.81%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.389"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -0.86%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.999"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -0.11%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.831"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +13.59%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'3.106.81"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +0.06%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.198"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -2.34%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  -2.55%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'93.792.78"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +1.87%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'34.83"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +1.10%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'  -1.50%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.686.80"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -0.12%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.097.21"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -1.60%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'3.66"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +0.54%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'14.86"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +0.72%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'5.98"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +3.17%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'445.69"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -0.10%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.0000201"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -0.48%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'8.98"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -3.36%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'8.19"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +4.06%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'5.63"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +0.08%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("B27").Value = "'Litecoin"
$ws.Range("B27").Style = "Normal"
$ws.Range("C27").Value = "'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("C27").Style = "Normal"
$ws.Range("D27").Value = "'86.20"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +6.77%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("B28").Value = "'Aptos"
$ws.Range("B28").Style = "Normal"
$ws.Range("C28").Value = "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("C28").Style = "Normal"
$ws.Range("D28").Value = "'12.17"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +4.74%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'3.271.42"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.22%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  +0.21%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'0.246"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +7.54%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'0.179"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +7.07%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.125"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -10.49%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'9.26"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -0.39%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  +0.06%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("B36").Value = "'Kaspa"
$ws.Range("B36").Style = "Normal"
$ws.Range("C36").Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("C36").Style = "Normal"
$ws.Range("D36").Value = "'0.163"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -3.37%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("B37").Value = "'RenderToken"
$ws.Range("B37").Style = "Normal"
$ws.Range("C37").Value = "'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("C37").Style = "Normal"
$ws.Range("D37").Value = "'7.88"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -1.10%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'26.13"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -0.37%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'1.91"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -1.26%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.453"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +4.87%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'3.81"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -7.26%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("B42").Value = "'WhiteBITCoin"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "'24.02"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +8.22%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("B43").Value = "'Bittensor"
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = "'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = "'476.86"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -0.58%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'1.28"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -1.41%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  -5.53%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D47").Value = "'160.52"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +1.04%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  -0.94%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  -2.78%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'4.46"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +2.18%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  -3.24%  "
$ws.Range("E51").Style = "Normal"
